$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.711.23'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '1.626.97'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.52'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.257'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0636'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('E10').Value = '  -4.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '1.654.29'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('D14').Value = '1.850.83'
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.552'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.59'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').Value = '25.720.35'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.00'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.42'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.78'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.92'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.20'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.30'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.32%  '
$ws.Range('E27').Value = '  -2.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.82'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.24'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0487'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.32'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.24'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.59'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.76%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.896'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.544'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.88%  '
$ws.Range('D39').Value = '1.106.46'
$ws.Range('E39').Value = '  -2.29%  '
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.08'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.796'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('D45').Value = '1.757.74'
$ws.Range('E45').Value = '  -1.38%  '
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.90'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.66'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.36'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.76%  '
$ws.Range('E51').Value = '  -0.74%  '
